$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 85

for ($r = 1; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)   # column A (source text) - stays as is
    $cCell = $ws.Cells.Item($r, 3)   # column C (old translation slot)
    $dCell = $ws.Cells.Item($r, 4)   # column D (old translation slot)
    $bCell = $ws.Cells.Item($r, 2)   # column B (new, consolidated translation slot)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    if ($cVal -ne $null -and $cVal -ne "") {
        $cCell.Copy($bCell)
    } elseif ($dVal -ne $null -and $dVal -ne "") {
        $dCell.Copy($bCell)
    } else {
        $aCell.Copy($bCell)
    }

    # Remove the now-obsolete columns C and D
    $cCell.Clear()
    $dCell.Clear()
}
